$d = $word.ActiveDocument

# Move to the very end of the document body (after the existing empty paragraph)
$r = $d.Content
$r.Collapse(0) | Out-Null

# Insert a fresh paragraph mark so the existing trailing empty <w:p/> is left untouched,
# then land the new range inside that freshly-created paragraph.
$r.InsertParagraphAfter()
$r.Collapse(0) | Out-Null
$r.Move(1, 1) | Out-Null

$startParaIndex = $d.Paragraphs.Count

$text = "[PUMP:TBV:1111]`rPUMP:HRD:3350`rDetails regarding the full color touchscreen. `rPUMP:HTP:1500`rTest 1500 `rPUMP:HRD:0000`rDetails regarding the size and weight of the pump. `r[PUMP:TBV:1]`rACE:SRS:1`rThe software shall provide a bolus feature which generates boluses in the range of 0.01 to 25 units, which an increment of 0.01 units.  `rPUMP:SDS:10`rHere are details of how the bolus calculator works ….         `rPUMP:SVAL:100`rThis test validates bolus features… blah, blah, blah        "
$r.InsertAfter($text)

# Apply paragraph-level formatting (style / left indent) to each newly inserted paragraph.
$specs = @(
    @{ Style = $null; Ind = $null },
    @{ Style = "ListBullet"; Ind = $null },
    @{ Style = $null; Ind = 360 },
    @{ Style = $null; Ind = 720 },
    @{ Style = $null; Ind = 720 },
    @{ Style = "ListBullet"; Ind = $null },
    @{ Style = $null; Ind = 360 },
    @{ Style = $null; Ind = $null },
    @{ Style = "ListBullet"; Ind = $null },
    @{ Style = $null; Ind = 360 },
    @{ Style = $null; Ind = 720 },
    @{ Style = $null; Ind = 720 },
    @{ Style = $null; Ind = 720 },
    @{ Style = $null; Ind = 720 }
)

for ($i = 0; $i -lt $specs.Count; $i++) {
    $p = $d.Paragraphs($startParaIndex + $i)
    $spec = $specs[$i]
    if ($spec.Style) {
        $p.Style = $spec.Style
    }
    if ($null -ne $spec.Ind) {
        $p.Range.ParagraphFormat.LeftIndent = [double]($spec.Ind) / 20.0
    }
}

Write-Host "Inserted paragraphs:" $specs.Count
Write-Host "Total paragraphs now:" $d.Paragraphs.Count
